# Build the "Passage" (race results) sheet: a 11-column header row plus
# 7 data rows, replacing the original tiny 2-row / 3-column sheet.
#
# NOTE on write order: new text is appended to the workbook's shared-string
# table in first-write order, so the cells below are populated strictly
# top-to-bottom / left-to-right (row 1, then row 2 col by col, etc.) to line
# up with how the source workbook was produced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header row --------------------------------------------------------
$ws.Range("A1").Value = "id_participants"
$ws.Range("B1").Value = "Nom"
$ws.Range("C1").Value = "Prenom"
$ws.Range("D1").Value = "id_epreuve"
$ws.Range("E1").Value = "Epreuve"
$ws.Range("F1").Value = "Date"
$ws.Range("G1").Value = "Temps 1"
$ws.Range("H1").Value = "Temps 2"
$ws.Range("I1").Value = "Meilleur Temps"
$ws.Range("J1").Value = "id_categorie"
$ws.Range("K1").Value = "Type"

# Helper-free, explicit rows (the runtime's PS subset has no function defs
# that reliably close over COM objects, so each row is spelled out).

function Set-DateText($cell) {
    # Force literal text storage ("2021-02-27" must stay a string, not get
    # auto-converted to a date serial) without leaving a custom number
    # format behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = "2021-02-27"
    $cell.ClearFormats()
}

# ---- row 2 --------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Génique"
$ws.Range("C2").Value = "Yoann"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = "La Descente des Alpes - M1"
Set-DateText $ws.Range("F2")
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "M1"

# ---- row 3 --------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Mairot"
$ws.Range("C3").Value = "Jean-christophe"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "La Descente des Alpes - M1"
Set-DateText $ws.Range("F3")
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = "M1"

# ---- row 4 --------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Cherief"
$ws.Range("C4").Value = "Saufiane"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "La Descente des Alpes - M1"
Set-DateText $ws.Range("F4")
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = "M1"

# ---- row 5 --------------------------------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Rameau"
$ws.Range("C5").Value = "Célia"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "La Descente des Alpes - M1"
Set-DateText $ws.Range("F5")
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = "M1"

# ---- row 6 --------------------------------------------------------------
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Ligourel"
$ws.Range("C6").Value = "Teedji"
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "La Descente des Alpes - M1"
Set-DateText $ws.Range("F6")
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = "M1"

# ---- row 7 --------------------------------------------------------------
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "toto"
$ws.Range("C7").Value = "tata"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "La Descente des Alpes - M1"
Set-DateText $ws.Range("F7")
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = "M1"

# ---- row 8 --------------------------------------------------------------
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "choula"
$ws.Range("C8").Value = "poula"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "La Descente des Alpes - M1"
Set-DateText $ws.Range("F8")
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = "M1"
